$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2173913043478261
$ws.Range("C2").Value = 0.5141776937618148
$ws.Range("J2").Value = 0.01701323251417769
$ws.Range("O2").Value = 0.001890359168241966
$ws.Range("P2").Value = 0.1474480151228733
$ws.Range("S2").Value = 0.1020793950850662

# Row 3
$ws.Range("B3").Value = 0.0176678445229682
$ws.Range("C3").Value = 0.04593639575971731
$ws.Range("J3").Value = 0.01060070671378092
$ws.Range("P3").Value = 0.7243816254416962
$ws.Range("S3").Value = 0.2014134275618374

# Row 4
$ws.Range("J4").Value = 0.02325581395348837
$ws.Range("P4").Value = 0.7093023255813954
$ws.Range("S4").Value = 0.2674418604651163

# Row 5
$ws.Range("J5").Value = 0.2
$ws.Range("P5").Value = 0.4
$ws.Range("S5").Value = 0.4

# Row 6
$ws.Range("B6").Value = 0.06304347826086956
$ws.Range("D6").Value = 0.008695652173913044
$ws.Range("E6").Value = 0.002173913043478261
$ws.Range("F6").Value = 0.09347826086956522
$ws.Range("J6").Value = 0.2260869565217391
$ws.Range("O6").Value = 0.02173913043478261
$ws.Range("Q6").Value = 0.1521739130434783
$ws.Range("R6").Value = 0.05217391304347826
$ws.Range("S6").Value = 0.3804347826086957

# Row 7
$ws.Range("B7").Value = 0.08551068883610451
$ws.Range("D7").Value = 0.02850356294536817
$ws.Range("E7").Value = 0.002375296912114014
$ws.Range("F7").Value = 0.09263657957244656
$ws.Range("J7").Value = 0.1163895486935867
$ws.Range("O7").Value = 0.02375296912114014
$ws.Range("Q7").Value = 0.1852731591448931
$ws.Range("R7").Value = 0.08076009501187649
$ws.Range("S7").Value = 0.3847980997624703

# Row 8
$ws.Range("B8").Value = 0.07821782178217822
$ws.Range("D8").Value = 0.01782178217821782
$ws.Range("F8").Value = 0.06336633663366337
$ws.Range("J8").Value = 0.1
$ws.Range("O8").Value = 0.03069306930693069
$ws.Range("Q8").Value = 0.1693069306930693
$ws.Range("R8").Value = 0.0891089108910891
$ws.Range("S8").Value = 0.4514851485148515

# Row 9
$ws.Range("B9").Value = 0.08823529411764706
$ws.Range("D9").Value = 0.01890756302521008
$ws.Range("F9").Value = 0.06512605042016807
$ws.Range("J9").Value = 0.0861344537815126
$ws.Range("O9").Value = 0.03571428571428571
$ws.Range("Q9").Value = 0.1785714285714286
$ws.Range("R9").Value = 0.09453781512605042
$ws.Range("S9").Value = 0.4327731092436975

# Row 10
$ws.Range("B10").Value = 0.09217171717171717
$ws.Range("D10").Value = 0.01978114478114478
$ws.Range("E10").Value = 0.001683501683501683
$ws.Range("F10").Value = 0.06355218855218855
$ws.Range("J10").Value = 0.1161616161616162
$ws.Range("O10").Value = 0.01220538720538721
$ws.Range("Q10").Value = 0.2213804713804714
$ws.Range("R10").Value = 0.07786195286195287
$ws.Range("S10").Value = 0.3952020202020202

# Row 11
$ws.Range("G11").Value = 0.1253687315634218
$ws.Range("J11").Value = 0.09587020648967552
$ws.Range("K11").Value = 0.1858407079646018
$ws.Range("L11").Value = 0.5693215339233039
$ws.Range("S11").Value = 0.02359882005899705

# Row 12
$ws.Range("G12").Value = 0.7248157248157249
$ws.Range("J12").Value = 0.1769041769041769
$ws.Range("K12").Value = 0.007371007371007371
$ws.Range("L12").Value = 0.03931203931203931
$ws.Range("S12").Value = 0.05159705159705159

# Row 13
$ws.Range("G13").Value = 0.6626506024096386
$ws.Range("J13").Value = 0.2891566265060241
$ws.Range("S13").Value = 0.04819277108433735

# Row 15
$ws.Range("F15").Value = 0.01735357917570499
$ws.Range("H15").Value = 0.1778741865509761
$ws.Range("I15").Value = 0.0737527114967462
$ws.Range("J15").Value = 0.2885032537960954
$ws.Range("K15").Value = 0.07809110629067245
$ws.Range("M15").Value = 0.01735357917570499
$ws.Range("N15").Value = 0.004338394793926247
$ws.Range("O15").Value = 0.08676789587852494
$ws.Range("S15").Value = 0.2559652928416486

# Row 16
$ws.Range("F16").Value = 0.02694610778443114
$ws.Range("H16").Value = 0.218562874251497
$ws.Range("I16").Value = 0.09880239520958084
$ws.Range("J16").Value = 0.3652694610778443
$ws.Range("K16").Value = 0.1317365269461078
$ws.Range("M16").Value = 0.01197604790419162
$ws.Range("O16").Value = 0.03892215568862276
$ws.Range("S16").Value = 0.1077844311377246

# Row 17
$ws.Range("F17").Value = 0.01837837837837838
$ws.Range("H17").Value = 0.2010810810810811
$ws.Range("I17").Value = 0.09513513513513513
$ws.Range("J17").Value = 0.3675675675675676
$ws.Range("K17").Value = 0.1156756756756757
$ws.Range("M17").Value = 0.01297297297297297
$ws.Range("N17").Value = 0.003243243243243243
$ws.Range("O17").Value = 0.07567567567567568
$ws.Range("S17").Value = 0.1102702702702703

# Row 18
$ws.Range("F18").Value = 0.01058201058201058
$ws.Range("H18").Value = 0.1693121693121693
$ws.Range("I18").Value = 0.1058201058201058
$ws.Range("J18").Value = 0.3597883597883598
$ws.Range("K18").Value = 0.1296296296296296
$ws.Range("M18").Value = 0.03174603174603174
$ws.Range("O18").Value = 0.0873015873015873
$ws.Range("S18").Value = 0.1058201058201058

# Row 19
$ws.Range("F19").Value = 0.01532710280373832
$ws.Range("H19").Value = 0.2261682242990654
$ws.Range("I19").Value = 0.1046728971962617
$ws.Range("J19").Value = 0.3502803738317757
$ws.Range("K19").Value = 0.1158878504672897
$ws.Range("M19").Value = 0.01906542056074766
$ws.Range("N19").Value = 0.0007476635514018691
$ws.Range("O19").Value = 0.05719626168224299
$ws.Range("S19").Value = 0.1106542056074766
